# Reorder the "units" section rows by staff & students:
# swap the Time/Unit/Classroom/Lecturer/Delivery Mode (columns B:F) between
# row pairs (3,6), (8,9) and (15,16). Column A (Day) is identical for both
# rows in every pair, so it is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r1 = $ws.Range("B3:F3")
$r2 = $ws.Range("B6:F6")
$v1 = $r1.Value2
$v2 = $r2.Value2
$r1.Value2 = $v2
$r2.Value2 = $v1

$r1 = $ws.Range("B8:F8")
$r2 = $ws.Range("B9:F9")
$v1 = $r1.Value2
$v2 = $r2.Value2
$r1.Value2 = $v2
$r2.Value2 = $v1

$r1 = $ws.Range("B15:F15")
$r2 = $ws.Range("B16:F16")
$v1 = $r1.Value2
$v2 = $r2.Value2
$r1.Value2 = $v2
$r2.Value2 = $v1
